$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new Kaspa buy recorded on 2025-09-12 as row 5.
# Column A stores the date as literal text (matching the existing rows,
# which are plain text like "08/29/2025", not real date values), so a
# leading apostrophe is used to force text entry and avoid Excel's
# automatic date parsing; the cell style is then reset to "Normal" so the
# quote-prefix formatting doesn't linger on the cell.
$ws.Range("A2").Copy()
$ws.Range("A5").PasteSpecial()
$ws.Range("A5").Value = "'09/12/2025"
$ws.Range("A5").Style = "Normal"

$ws.Range("B5").Value = 862.7109999999993
$ws.Range("C5").Value = 0.05795683606677095
$ws.Range("D5").Value = 25
